$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes existing rows 3-21 down to 4-22),
# copying the formatting of the row above so the date style (column D)
# carries over correctly.
$ws.Rows.Item(3).Insert()

# Mercado ID
$ws.Cells.Item(3, 1).Value = 11
# Mercado
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
# Región
$ws.Cells.Item(3, 3).Value = "Bíobío"
# Fecha
$ws.Cells.Item(3, 4).Value = Get-Date -Year 2021 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Codreg
$ws.Cells.Item(3, 5).Value = 8
# Categoría ID
$ws.Cells.Item(3, 6).Value = 100112030
# Categoría
$ws.Cells.Item(3, 7).Value = "Poroto granado"
# Variedad
$ws.Cells.Item(3, 8).Value = "Sin especificar"
# Calidad
$ws.Cells.Item(3, 9).Value = "Primera"
# Volumen
$ws.Cells.Item(3, 10).Value = 100
# Precio mínimo
$ws.Cells.Item(3, 11).Value = 38000
# Precio máximo
$ws.Cells.Item(3, 12).Value = 40000
# Precio promedio ponderado
$ws.Cells.Item(3, 13).Value = 39000
# Unidad de comercialización
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
# Origen
$ws.Cells.Item(3, 15).Value = "Provincia de Limarí"
# Precio $/Kg
$ws.Cells.Item(3, 16).Value = 1560
# Kg o Unidades
$ws.Cells.Item(3, 17).Value = 25
# Clasificación
$ws.Cells.Item(3, 18).Value = "Hortaliza"
